# Append a new bullet item to the end of the journal list:
#   "Attended 2nd weekly meeting."
# with "nd" rendered as a superscript, matching the style/numbering
# of the preceding list paragraphs.

$d = $word.ActiveDocument

# Move to the very end of the document and add a new paragraph there.
# InsertParagraphAfter on the end-of-story range inherits the paragraph
# style (ListParagraph), numbering (ilvl 0 / numId 1) and run formatting
# (sz/szCs 24) from the last existing paragraph, which is exactly what we
# want for the new bullet.
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last

# Insert the full sentence as plain text first (inherits the paragraph's
# default run formatting: sz=24 / szCs=24, no superscript).
$r = $newPara.Range
$r.Collapse(0)
$r.InsertAfter("Attended 2nd weekly meeting.")

# Now mark just the "nd" ordinal suffix as superscript. Locate it by
# character offset relative to the paragraph start ("Attended 2" is 10
# characters, followed by the 2-character "nd").
$paraStart = $newPara.Range.Start
$ndStart = $paraStart + 10
$ndEnd = $paraStart + 12
$ndRange = $d.Range($ndStart, $ndEnd)
$ndRange.Font.Superscript = $true
